# Weekly update: insert a new pair of rows (Primera / Segunda) for Betarraga
# at Terminal La Palmera de La Serena, ahead of the existing rows starting
# at row 320, shifting the old rows 320-339 down to 322-341.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the 320 block; everything
# currently at row 320 and below shifts down by two rows.
$ws.Range("A320:A321").EntireRow.Insert()

# Common (constant-across-the-table) values for this market/category block.
$mercadoId = 8
$mercado   = "Terminal La Palmera de La Serena"
$region    = "Coquimbo"
$codreg    = 4
$categoriaId = 100114014
$categoria = "Betarraga"
$variedad  = "Sin especificar"
$unidad    = "`$/paquete 3 unidades"
$origen    = "Provincia del Elquí"
$kgUnidades = 3
$clasificacion = "Hortaliza"

# New row 320: Primera quality
$ws.Cells.Item(320, 1).Value = $mercadoId
$ws.Cells.Item(320, 2).Value = $mercado
$ws.Cells.Item(320, 3).Value = $region
$ws.Cells.Item(320, 4).Value = 44826
$ws.Cells.Item(320, 5).Value = $codreg
$ws.Cells.Item(320, 6).Value = $categoriaId
$ws.Cells.Item(320, 7).Value = $categoria
$ws.Cells.Item(320, 8).Value = $variedad
$ws.Cells.Item(320, 9).Value = "Primera"
$ws.Cells.Item(320, 10).Value = 2400
$ws.Cells.Item(320, 11).Value = 550
$ws.Cells.Item(320, 12).Value = 600
$ws.Cells.Item(320, 13).Value = 575
$ws.Cells.Item(320, 14).Value = $unidad
$ws.Cells.Item(320, 15).Value = $origen
$ws.Cells.Item(320, 16).Value = 192
$ws.Cells.Item(320, 17).Value = $kgUnidades
$ws.Cells.Item(320, 18).Value = $clasificacion

# New row 321: Segunda quality
$ws.Cells.Item(321, 1).Value = $mercadoId
$ws.Cells.Item(321, 2).Value = $mercado
$ws.Cells.Item(321, 3).Value = $region
$ws.Cells.Item(321, 4).Value = 44826
$ws.Cells.Item(321, 5).Value = $codreg
$ws.Cells.Item(321, 6).Value = $categoriaId
$ws.Cells.Item(321, 7).Value = $categoria
$ws.Cells.Item(321, 8).Value = $variedad
$ws.Cells.Item(321, 9).Value = "Segunda"
$ws.Cells.Item(321, 10).Value = 1600
$ws.Cells.Item(321, 11).Value = 450
$ws.Cells.Item(321, 12).Value = 500
$ws.Cells.Item(321, 13).Value = 475
$ws.Cells.Item(321, 14).Value = $unidad
$ws.Cells.Item(321, 15).Value = $origen
$ws.Cells.Item(321, 16).Value = 158
$ws.Cells.Item(321, 17).Value = $kgUnidades
$ws.Cells.Item(321, 18).Value = $clasificacion
